# Update the "想去人数" (F column) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 3429
    3  = 154
    5  = 1764
    7  = 484
    8  = 389
    14 = 11
    15 = 63
    22 = 138
    23 = 30
    24 = 422
    25 = 301
    26 = 124
    28 = 21
    30 = 604
    31 = 2410
    34 = 494
    35 = 694
    41 = 560
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
